$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) - reorder labels
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "bedrooms_2"
$ws.Range("D1").Value = "kitchens_2"
$ws.Range("E1").Value = "living_rooms_1"
$ws.Range("F1").Value = "living_rooms_2"

# Update data rows 2-7 with the new 0/1 matrix
$data = @(
    @(0,1,0,0,0,0),
    @(1,0,0,0,0,0),
    @(0,0,0,0,0,1),
    @(0,0,1,0,0,0),
    @(0,0,0,1,0,0),
    @(0,0,0,0,1,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowValues[$col - 1]
    }
}
